# Add 2022-Q1 data:
#   1. Insert a new "2022-Q1" worksheet (holdings detail) right before "总计".
#   2. Insert the corresponding aggregate row at the top of "总计" (and renumber its index column).

$wb = $excel.ActiveWorkbook

function Set-TextCell($rng, [string]$text) {
    # Force the cell to store `text` as TEXT even when it looks numeric
    # (e.g. "515400", "2.29"), matching the source data's inlineStr cells.
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

# ---------------------------------------------------------------------------
# 1. New "2022-Q1" worksheet
# ---------------------------------------------------------------------------

# Duplicate an existing quarter sheet so the new sheet starts with the exact
# same column headers / fonts / borders already used by its siblings, then
# drop it immediately in front of "总计".
$wsTemplate = $wb.Worksheets.Item("2021-Q4")
$zjSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsTemplate.Copy($zjSheet)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# The template only has 2 data rows; this sheet needs 3, so clone row 3's
# formatting (an unstyled data row, except the bordered index cell in col A)
# down into row 4.
$newSheet.Range("A3:H3").Copy()
$newSheet.Range("A4:H4").PasteSpecial(-4122)

# Row 2
Set-TextCell $newSheet.Range("B2") "515400"
Set-TextCell $newSheet.Range("C2") "富国中证大数据产业ETF"
Set-TextCell $newSheet.Range("D2") "2.29"
Set-TextCell $newSheet.Range("E2") "99.33"
Set-TextCell $newSheet.Range("F2") "2.54"
Set-TextCell $newSheet.Range("G2") "0.0582"
$newSheet.Range("H2").Value = 10

# Row 3
Set-TextCell $newSheet.Range("B3") "560660"
Set-TextCell $newSheet.Range("C3") "新华中证云计算50交易型开放式指数证券投资基金"
Set-TextCell $newSheet.Range("D3") "1.75"
Set-TextCell $newSheet.Range("E3") "97.21"
Set-TextCell $newSheet.Range("F3") "2.89"
Set-TextCell $newSheet.Range("G3") "0.0506"
$newSheet.Range("H3").Value = 9

# Row 4
$newSheet.Range("A4").Value = 2
Set-TextCell $newSheet.Range("B4") "516000"
Set-TextCell $newSheet.Range("C4") "华夏中证大数据产业ETF"
Set-TextCell $newSheet.Range("D4") "0.58"
Set-TextCell $newSheet.Range("E4") "97.61"
Set-TextCell $newSheet.Range("F4") "2.51"
Set-TextCell $newSheet.Range("G4") "0.0146"
$newSheet.Range("H4").Value = 10

# ---------------------------------------------------------------------------
# 2. Update "总计" with the new 2022-Q1 aggregate row
# ---------------------------------------------------------------------------

$zj = $wb.Worksheets.Item("总计")

# Remember the existing 4 data rows before they get shifted down.
$bVals = @($zj.Range("B2").Value2, $zj.Range("B3").Value2, $zj.Range("B4").Value2, $zj.Range("B5").Value2)
$cVals = @($zj.Range("C2").Value2, $zj.Range("C3").Value2, $zj.Range("C4").Value2, $zj.Range("C5").Value2)
$dVals = @($zj.Range("D2").Value2, $zj.Range("D3").Value2, $zj.Range("D4").Value2, $zj.Range("D5").Value2)

# Extend the formatting of the index column (styled, bordered) down to the
# new last row (row 6).
$zj.Range("A5:D5").Copy()
$zj.Range("A6:D6").PasteSpecial(-4122)

# Shift the old rows 2..5 down to rows 3..6.
for ($i = 0; $i -lt 4; $i++) {
    $destRow = $i + 3
    $zj.Range("B$destRow").Value2 = $bVals[$i]
    $zj.Range("C$destRow").Value2 = $cVals[$i]
    $zj.Range("D$destRow").Value2 = $dVals[$i]
}

# New 2022-Q1 summary row at the top (row 2).
$zj.Range("B2").Value2 = "2022-Q1"
$zj.Range("C2").Value2 = 3
$zj.Range("D2").Value2 = 0.12

# Renumber the 0-based index column for all 5 data rows.
for ($i = 0; $i -lt 5; $i++) {
    $zj.Range("A" + ($i + 2)).Value2 = $i
}
